# Build v2.1.2: Fix SearchCriteria variants and Schemas sheet grouping/sorting
#
# Collapses the inline request/response body field listings on the
# "Body", "200", "204", "400" sheets down to a single schema-reference
# row, and adds a matching schema-reference row (errorResponse1) on the
# plain error-response sheets ("401","403","404","429","500") that
# previously only had header rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# "Body" sheet: request body -> single setCalendar.211207Request ref
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Body")
$ws.Range("B3").Value2 = "setCalendar.211207Request"
$ws.Range("D3").Value2 = ""
$ws.Range("E3").Value2 = "schema"
$ws.Range("G3").Value2 = "setCalendar.211207Request"
$ws.Range("I3").Value2 = "Yes"
$ws.Range("L3").Value2 = ""
$ws.Range("O3").Value2 = ""
$ws.Rows("4:11").Delete()

# ---------------------------------------------------------------
# "200" sheet: response body -> single setCalendar.211207Response ref
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("200")
$ws.Range("B3").Value2 = "setCalendar.211207Response"
$ws.Range("D3").Value2 = ""
$ws.Range("E3").Value2 = "schema"
$ws.Range("G3").Value2 = "setCalendar.211207Response"
$ws.Range("I3").Value2 = "Yes"
$ws.Range("L3").Value2 = ""
$ws.Range("O3").Value2 = ""
$ws.Rows("4:4").Delete()

# ---------------------------------------------------------------
# "204" sheet: add the setCalendar.211207Response schema ref row
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("204")
$ws.Range("A3").Value2 = "content"
$ws.Range("B3").Value2 = "setCalendar.211207Response"
$ws.Range("E3").Value2 = "schema"
$ws.Range("G3").Value2 = "setCalendar.211207Response"
$ws.Range("I3").Value2 = "Yes"

# ---------------------------------------------------------------
# "400" sheet: response body -> single errorResponse ref
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("400")
$ws.Range("B3").Value2 = "errorResponse"
$ws.Range("D3").Value2 = ""
$ws.Range("E3").Value2 = "schema"
$ws.Range("G3").Value2 = "errorResponse"
$ws.Range("I3").Value2 = "Yes"
$ws.Range("L3").Value2 = ""
$ws.Range("O3").Value2 = ""
$ws.Rows("4:6").Delete()

# ---------------------------------------------------------------
# "401","403","404","429","500": add the errorResponse1 schema ref row
# ---------------------------------------------------------------
$errorSheets = @("401", "403", "404", "429", "500")
foreach ($name in $errorSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A3").Value2 = "content"
    $ws.Range("B3").Value2 = "errorResponse1"
    $ws.Range("E3").Value2 = "schema"
    $ws.Range("G3").Value2 = "errorResponse1"
    $ws.Range("I3").Value2 = "Yes"
}
